$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Shared-string text updates (order matters for shared-string table
#    layout: modify the existing "more sfx" text in place first, then
#    introduce the two brand-new strings in the order they should be
#    appended).
# ---------------------------------------------------------------------

# M7 (and M13, which currently also points at the same shared string)
# get the text updated in place first.
$ws.Range("M7").Value = "more sfx, more animations"

# M13 now becomes its own, distinct string.
$ws.Range("M13").Value = "more sfx, outside, more animations"

# New cell M5 with a brand-new string.
$ws.Range("M5").Value = "sound effects, models, details"

# ---------------------------------------------------------------------
# 2) New / updated numeric cells in column L (%DoneSat for RL Tasks)
# ---------------------------------------------------------------------
$ws.Range("L5").Value = 0.85
$ws.Range("L5").NumberFormat = "0%"

$ws.Range("L6").Value = 0.85
$ws.Range("L6").NumberFormat = "0%"

$ws.Range("L7").Value = 0.95
$ws.Range("L8").Value = 0.95
$ws.Range("L9").Value = 0.95

$ws.Range("L10").Value = 0
$ws.Range("L10").NumberFormat = "0%"

$ws.Range("L11").Value = 0
$ws.Range("L11").NumberFormat = "0%"

$ws.Range("L12").Value = 0
$ws.Range("L12").NumberFormat = "0%"

$ws.Range("L15").Value = 0
$ws.Range("L15").NumberFormat = "0%"

$ws.Range("L16").Value = 0
$ws.Range("L16").NumberFormat = "0%"

# ---------------------------------------------------------------------
# 3) New text cell K11 ("Change in enviroment" for task T8)
# ---------------------------------------------------------------------
$ws.Range("K11").Value = "H9"

# ---------------------------------------------------------------------
# 4) Strikethrough formatting for the "Other" block (task O1 done) and
#    for RL Task T2 (row 6) / Horror-Element H2 (row 5).
#    Apply the plain (borderless) strikethrough cells first, then the
#    bordered ones, so the newly created styles line up with the
#    expected style indices (16 = plain strike, 17 = bordered strike).
# ---------------------------------------------------------------------
foreach ($addr in @("U4", "V4", "O5")) {
    $ws.Range($addr).Font.Strikethrough = $true
}

foreach ($addr in @("T4", "N5", "G6")) {
    $ws.Range($addr).Font.Strikethrough = $true
    $ws.Range($addr).Borders.Item(7).LineStyle = 1
}

# ---------------------------------------------------------------------
# 5) Update the active selection to match the saved view state.
# ---------------------------------------------------------------------
$null = $ws.Range("G6").Select()
